$d = $word.ActiveDocument

$d.Content.Find.Execute('548÷5=109, 3', $true, $false, $false, $false, $false, $true, 1, $false, '216÷3=72, 0', 2) | Out-Null
$d.Content.Find.Execute('693÷7=99, 0', $true, $false, $false, $false, $false, $true, 1, $false, '640÷4=160, 0', 2) | Out-Null
$d.Content.Find.Execute('216÷9=24, 0', $true, $false, $false, $false, $false, $true, 1, $false, '219÷3=73, 0', 2) | Out-Null
$d.Content.Find.Execute('248÷7=35, 3', $true, $false, $false, $false, $false, $true, 1, $false, '593÷6=98, 5', 2) | Out-Null
$d.Content.Find.Execute('510÷8=63, 6', $true, $false, $false, $false, $false, $true, 1, $false, '759÷4=189, 3', 2) | Out-Null
$d.Content.Find.Execute('477÷6=79, 3', $true, $false, $false, $false, $false, $true, 1, $false, '314÷5=62, 4', 2) | Out-Null
$d.Content.Find.Execute('385÷3=128, 1', $true, $false, $false, $false, $false, $true, 1, $false, '227÷3=75, 2', 2) | Out-Null
$d.Content.Find.Execute('279÷2=139, 1', $true, $false, $false, $false, $false, $true, 1, $false, '207÷8=25, 7', 2) | Out-Null
$d.Content.Find.Execute('702÷4=175, 2', $true, $false, $false, $false, $false, $true, 1, $false, '284÷4=71, 0', 2) | Out-Null
$d.Content.Find.Execute('921÷4=230, 1', $true, $false, $false, $false, $false, $true, 1, $false, '477÷4=119, 1', 2) | Out-Null
$d.Content.Find.Execute('352÷8=44, 0', $true, $false, $false, $false, $false, $true, 1, $false, '365÷3=121, 2', 2) | Out-Null
$d.Content.Find.Execute('739÷8=92, 3', $true, $false, $false, $false, $false, $true, 1, $false, '988÷8=123, 4', 2) | Out-Null
$d.Content.Find.Execute('599÷5=119, 4', $true, $false, $false, $false, $false, $true, 1, $false, '329÷9=36, 5', 2) | Out-Null
$d.Content.Find.Execute('786÷8=98, 2', $true, $false, $false, $false, $false, $true, 1, $false, '158÷9=17, 5', 2) | Out-Null
$d.Content.Find.Execute('767÷6=127, 5', $true, $false, $false, $false, $false, $true, 1, $false, '215÷7=30, 5', 2) | Out-Null
$d.Content.Find.Execute('440÷4=110, 0', $true, $false, $false, $false, $false, $true, 1, $false, '241÷4=60, 1', 2) | Out-Null
$d.Content.Find.Execute('195÷9=21, 6', $true, $false, $false, $false, $false, $true, 1, $false, '662÷2=331, 0', 2) | Out-Null
$d.Content.Find.Execute('533÷4=133, 1', $true, $false, $false, $false, $false, $true, 1, $false, '605÷4=151, 1', 2) | Out-Null
$d.Content.Find.Execute('247÷8=30, 7', $true, $false, $false, $false, $false, $true, 1, $false, '143÷7=20, 3', 2) | Out-Null
$d.Content.Find.Execute('546÷3=182, 0', $true, $false, $false, $false, $false, $true, 1, $false, '199÷8=24, 7', 2) | Out-Null
$d.Content.Find.Execute('718÷5=143, 3', $true, $false, $false, $false, $false, $true, 1, $false, '731÷6=121, 5', 2) | Out-Null
$d.Content.Find.Execute('801÷7=114, 3', $true, $false, $false, $false, $false, $true, 1, $false, '325÷4=81, 1', 2) | Out-Null
$d.Content.Find.Execute('986÷9=109, 5', $true, $false, $false, $false, $false, $true, 1, $false, '586÷7=83, 5', 2) | Out-Null
$d.Content.Find.Execute('318÷5=63, 3', $true, $false, $false, $false, $false, $true, 1, $false, '910÷8=113, 6', 2) | Out-Null
$d.Content.Find.Execute('302÷2=151, 0', $true, $false, $false, $false, $false, $true, 1, $false, '811÷9=90, 1', 2) | Out-Null
